$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.74669534304606
$ws.Range("C2").Value = 9.570545766966038
$ws.Range("D2").Value = 8.812900301709805
$ws.Range("F2").Value = 34.32704813639048
$ws.Range("G2").Value = 3.667690417862037
$ws.Range("I2").Value = 25.55654061914504
$ws.Range("J2").Value = 10.45043794749669
$ws.Range("L2").Value = 11.84404138942265
$ws.Range("M2").Value = 16.61903971452296
$ws.Range("O2").Value = 26.07417623181908
$ws.Range("B3").Value = 16.23557823278766
$ws.Range("C3").Value = 9.280665729017587
$ws.Range("D3").Value = 8.818693701607035
$ws.Range("F3").Value = 34.45886348785208
$ws.Range("G3").Value = 3.66976762776137
$ws.Range("I3").Value = 25.702039565122
$ws.Range("J3").Value = 10.47237552372833
$ws.Range("L3").Value = 11.84538184591604
$ws.Range("M3").Value = 16.50074004854826
$ws.Range("O3").Value = 26.18665390332444
$ws.Range("B4").Value = 15.91445091172337
$ws.Range("C4").Value = 9.096975874180117
$ws.Range("D4").Value = 8.823065909620675
$ws.Range("F4").Value = 34.54830920163663
$ws.Range("G4").Value = 3.671111111159346
$ws.Range("I4").Value = 25.7967410514179
$ws.Range("J4").Value = 10.48655737046804
$ws.Range("L4").Value = 11.84746529509748
$ws.Range("M4").Value = 16.42926177312293
$ws.Range("O4").Value = 26.26190201513068
$ws.Range("B5").Value = 15.78193647569507
$ws.Range("C5").Value = 9.020766792642592
$ws.Range("D5").Value = 8.825053089633277
$ws.Range("F5").Value = 34.5868934610488
$ws.Range("G5").Value = 3.671675761637843
$ws.Range("I5").Value = 25.83668228450528
$ws.Range("J5").Value = 10.49251618043143
$ws.Range("L5").Value = 11.84863227825399
$ws.Range("M5").Value = 16.40044623818728
$ws.Range("O5").Value = 26.29411836295638
$ws.Range("B6").Value = 15.75983837796769
$ws.Range("C6").Value = 9.008033029663236
$ws.Range("D6").Value = 8.825395484099756
$ws.Range("F6").Value = 34.59342909244285
$ws.Range("G6").Value = 3.671770560033802
$ws.Range("I6").Value = 25.84339602925295
$ws.Range("J6").Value = 10.49351649889046
$ws.Range("L6").Value = 11.84884528753766
$ws.Range("M6").Value = 16.39568092116727
$ws.Range("O6").Value = 26.29956150547707
$ws.Range("B7").Value = 15.91267021411448
$ws.Range("C7").Value = 9.095953463275203
$ws.Range("D7").Value = 8.823091876802735
$ws.Range("F7").Value = 34.54882092793341
$ws.Range("G7").Value = 3.671118656641756
$ws.Range("I7").Value = 25.79727424744836
$ws.Range("J7").Value = 10.48663700519399
$ws.Range("L7").Value = 11.84747974472638
$ws.Range("M7").Value = 16.42887186403903
$ws.Range("O7").Value = 26.26233021618858
$ws.Range("B8").Value = 16.57208254107206
$ws.Range("C8").Value = 9.47182443050357
$ws.Range("D8").Value = 8.814729023267514
$ws.Range("F8").Value = 34.37072836369916
$ws.Range("G8").Value = 3.668392544736989
$ws.Range("I8").Value = 25.60559562747599
$ws.Range("J8").Value = 10.45785455134427
$ws.Range("L8").Value = 11.84424261235128
$ws.Range("M8").Value = 16.57802246641199
$ws.Range("O8").Value = 26.11167247334984
$ws.Range("B9").Value = 17.7999577297981
$ws.Range("C9").Value = 10.16050416138412
$ws.Range("D9").Value = 8.804773302631066
$ws.Range("F9").Value = 34.08926224068957
$ws.Range("G9").Value = 3.663584281851104
$ws.Range("I9").Value = 25.27225001739712
$ws.Range("J9").Value = 10.40703798346343
$ws.Range("L9").Value = 11.84784897430152
$ws.Range("M9").Value = 16.87874618471635
$ws.Range("O9").Value = 25.86547311283117
$ws.Range("B10").Value = 18.65353284127874
$ws.Range("C10").Value = 10.63322195483403
$ws.Range("D10").Value = 8.8013557316368
$ws.Range("F10").Value = 33.9241102485711
$ws.Range("G10").Value = 3.660375968918055
$ws.Range("I10").Value = 25.05323054692492
$ws.Range("J10").Value = 10.37309849228092
$ws.Range("L10").Value = 11.8565025710521
$ws.Range("M10").Value = 17.10341764703665
$ws.Range("O10").Value = 25.7148165718911
$ws.Range("B11").Value = 19.02973062302632
$ws.Range("C11").Value = 10.84038871752959
$ws.Range("D11").Value = 8.800640058917329
$ws.Range("F11").Value = 33.85808578905639
$ws.Range("G11").Value = 3.65898611538062
$ws.Range("I11").Value = 24.95920736216943
$ws.Range("J11").Value = 10.35838856768797
$ws.Range("L11").Value = 11.86172793799511
$ws.Range("M11").Value = 17.20615626591039
$ws.Range("O11").Value = 25.65288612556399
$ws.Range("B12").Value = 19.17033387290493
$ws.Range("C12").Value = 10.91765744431055
$ws.Range("D12").Value = 8.800489048301234
$ws.Range("F12").Value = 33.83439824451109
$ws.Range("G12").Value = 3.658469770448328
$ws.Range("I12").Value = 24.92440971318452
$ws.Range("J12").Value = 10.3529226313584
$ws.Range("L12").Value = 11.86389054903701
$ws.Range("M12").Value = 17.24511374754463
$ws.Range("O12").Value = 25.63038791843996
$ws.Range("B13").Value = 19.14013663447351
$ws.Range("C13").Value = 10.90106947228777
$ws.Range("D13").Value = 8.800516244735231
$ws.Range("F13").Value = 33.83944123677488
$ws.Range("G13").Value = 3.658580532269159
$ws.Range("I13").Value = 24.93186811868797
$ws.Range("J13").Value = 10.35409518327985
$ws.Range("L13").Value = 11.86341663792606
$ws.Range("M13").Value = 17.23672160354556
$ws.Range("O13").Value = 25.63519084148687
$ws.Range("B14").Value = 19.0413359223137
$ws.Range("C14").Value = 10.84676958572232
$ws.Range("D14").Value = 8.800625234022718
$ws.Range("F14").Value = 33.85611062662781
$ws.Range("G14").Value = 3.658943435992144
$ws.Range("I14").Value = 24.9563283658948
$ws.Range("J14").Value = 10.35793679237556
$ws.Range("L14").Value = 11.86190218069228
$ws.Range("M14").Value = 17.20936040611713
$ws.Range("O14").Value = 25.6510160528836
$ws.Range("B15").Value = 18.98057289753735
$ws.Range("C15").Value = 10.81335421477331
$ws.Range("D15").Value = 8.800707600875892
$ws.Range("F15").Value = 33.86649244797167
$ws.Range("G15").Value = 3.659167020944745
$ws.Range("I15").Value = 24.97141606200146
$ws.Range("J15").Value = 10.36030346966027
$ws.Range("L15").Value = 11.86099843392327
$ws.Range("M15").Value = 17.19260701268447
$ws.Range("O15").Value = 25.66083374039284
$ws.Range("B16").Value = 18.62869387179595
$ws.Range("C16").Value = 10.61952039016885
$ws.Range("D16").Value = 8.801419330523386
$ws.Range("F16").Value = 33.92860876070097
$ws.Range("G16").Value = 3.660468195896369
$ws.Range("I16").Value = 25.05948808292175
$ws.Range("J16").Value = 10.37407445488764
$ws.Range("L16").Value = 11.85618689278521
$ws.Range("M16").Value = 17.09671211746148
$ws.Range("O16").Value = 25.7189970929487
$ws.Range("B17").Value = 18.4096430786638
$ws.Range("C17").Value = 10.49855643472638
$ws.Range("D17").Value = 8.802070388254526
$ws.Range("F17").Value = 33.96905074427385
$ws.Range("G17").Value = 3.661284221699868
$ws.Range("I17").Value = 25.11495434924682
$ws.Range("J17").Value = 10.38270895995884
$ws.Range("L17").Value = 11.85356425346605
$ws.Range("M17").Value = 17.03800333275372
$ws.Range("O17").Value = 25.75637263239559
$ws.Range("B18").Value = 18.28251820226273
$ws.Range("C18").Value = 10.42824247016253
$ws.Range("D18").Value = 8.80252387023385
$ws.Range("F18").Value = 33.99316851282956
$ws.Range("G18").Value = 3.66176013443887
$ws.Range("I18").Value = 25.14738508854926
$ws.Range("J18").Value = 10.38774397452395
$ws.Range("L18").Value = 11.85217718041953
$ws.Range("M18").Value = 17.00428735081974
$ws.Range("O18").Value = 25.77849142372595
$ws.Range("B19").Value = 18.23928508094291
$ws.Range("C19").Value = 10.40431006448422
$ws.Range("D19").Value = 8.802691000929217
$ws.Range("F19").Value = 34.00148133137047
$ws.Range("G19").Value = 3.661922397875085
$ws.Range("I19").Value = 25.15845626136849
$ws.Range("J19").Value = 10.38946055382391
$ws.Range("L19").Value = 11.85172843479051
$ws.Range("M19").Value = 16.99288132725629
$ws.Range("O19").Value = 25.78608706840403
$ws.Range("B20").Value = 18.43307950802142
$ws.Range("C20").Value = 10.51151007488251
$ws.Range("D20").Value = 8.801992909345062
$ws.Range("F20").Value = 33.96465693069915
$ws.Range("G20").Value = 3.661196676169347
$ws.Range("I20").Value = 25.10899522229736
$ws.Range("J20").Value = 10.38178269799194
$ws.Range("L20").Value = 11.85383088375056
$ws.Range("M20").Value = 17.04424780205112
$ws.Range("O20").Value = 25.75232960475845
$ws.Range("B21").Value = 19.07040727709229
$ws.Range("C21").Value = 10.86275118486519
$ws.Range("D21").Value = 8.800589969775324
$ws.Range("F21").Value = 33.8511787067645
$ws.Range("G21").Value = 3.658836572386118
$ws.Range("I21").Value = 24.94912189755004
$ws.Range("J21").Value = 10.3568055892928
$ws.Range("L21").Value = 11.86234203502812
$ws.Range("M21").Value = 17.21739581879288
$ws.Range("O21").Value = 25.64634189625536
$ws.Range("B22").Value = 19.47608255107175
$ws.Range("C22").Value = 11.08540608471197
$ws.Range("D22").Value = 8.800372195700756
$ws.Range("F22").Value = 33.78467824025021
$ws.Range("G22").Value = 3.657352154186762
$ws.Range("I22").Value = 24.84933877142337
$ws.Range("J22").Value = 10.3410898656688
$ws.Range("L22").Value = 11.86897562018373
$ws.Range("M22").Value = 17.33085341164186
$ws.Range("O22").Value = 25.58263232555674
$ws.Range("B23").Value = 19.26059410237652
$ws.Range("C23").Value = 10.96721706288887
$ws.Range("D23").Value = 8.800424676161263
$ws.Range("F23").Value = 33.81946782930871
$ws.Range("G23").Value = 3.658139121127475
$ws.Range("I23").Value = 24.90216445889192
$ws.Range("J23").Value = 10.34942214549398
$ws.Range("L23").Value = 11.86533764781176
$ws.Range("M23").Value = 17.27027993759666
$ws.Range("O23").Value = 25.61612536618838
$ws.Range("B24").Value = 18.42248759496351
$ws.Range("C24").Value = 10.50565612940271
$ws.Range("D24").Value = 8.802027690886632
$ws.Range("F24").Value = 33.96664067313659
$ws.Range("G24").Value = 3.661236234416097
$ws.Range("I24").Value = 25.11168765389308
$ws.Range("J24").Value = 10.38220124006927
$ws.Range("L24").Value = 11.85370996407082
$ws.Range("M24").Value = 17.04142456199396
$ws.Range("O24").Value = 25.75415549199552
$ws.Range("B25").Value = 17.47573002809403
$ws.Range("C25").Value = 9.979792256137932
$ws.Range("D25").Value = 8.806779884748037
$ws.Range("F25").Value = 34.15811634962577
$ws.Range("G25").Value = 3.664827843781863
$ws.Range("I25").Value = 25.35787969901158
$ws.Range("J25").Value = 10.42018648518033
$ws.Range("L25").Value = 11.84581403797063
$ws.Range("M25").Value = 16.7966427137085
$ws.Range("O25").Value = 25.9267834171541
